$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new Price (column D) value, only for rows whose price
# actually changed in this refresh. Values are plain text in the sheet
# (column D holds numeric-looking strings, not real numbers), so each one
# is written with a leading apostrophe to force a text entry instead of
# letting Excel auto-convert it to a number.
$priceUpdates = @{
    2  = "265.04"
    3  = "22.69"
    4  = "6.286"
    5  = "0.06136"
    6  = "3.600"
    7  = "6.670"
    8  = "1.346"
    9  = "0.8261"
    10 = "0.01354"
    11 = "0.1593"
    12 = "0.08223"
    13 = "0.03411"
    14 = "0.03131"
    15 = "0.09254"
    17 = "0.001708"
    18 = "0.04877"
    19 = "0.006235"
    20 = "0.005273"
    21 = "0.001088"
    23 = "3.766"
    24 = "2.289"
    25 = "0.3376"
    26 = "0.1238"
    27 = "0.0002680"
    40 = "0.04603"
    41 = "0.007014"
    42 = "0.1137"
    43 = "0.003399"
    44 = "0.01080"
    45 = "0.00006150"
    47 = "0.7782"
    48 = "0.1941"
    49 = "0.00002099"
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
}

# Every data row (2 through 51) also has its "Hora" (column G) value
# updated from "4" to "5".
for ($row = 2; $row -le 51; $row++) {
    $ws.Range("G$row").Value = "'5"
}
